# "added battery and button" -- append two component rows under the
# existing header (Component | Link | Mech/Elec/Firm).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Button"
$ws.Range("B2").Value = "https://www.digikey.ca/en/products/detail/adam-tech/SW-PB2-2EZ-A-RR-3-L1/15284423"
$ws.Range("C2").Value = "Elec"

$ws.Range("A3").Value = "Battery"
$ws.Range("B3").Value = "https://www.18650batterystore.com/en-ca/products/ydl-14500-battery"
$ws.Range("C3").Value = "Elec"

# Resize columns to fit the new (much longer) link text, same as the
# author's sheet where columns A-C are all bestFit/customWidth.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("B7").Select() | Out-Null
